# Apply scheduled price/profit refresh to all sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 14301.667
$ws.Range("I33").Value = 16417.309
$ws.Range("K33").Value = 16417.309
$ws.Range("M33").Value = -16188.309
# Row 40
$ws.Range("H40").Value = 3341467
$ws.Range("I40").Value = 3341467
$ws.Range("K40").Value = 3341467
$ws.Range("M40").Value = -3341292
# Row 138
$ws.Range("H138").Value = 13584.846
$ws.Range("I138").Value = 3334.3333
$ws.Range("J138").Value = 16660
$ws.Range("K138").Value = 10002.9999
$ws.Range("L138").Value = 49980
$ws.Range("M138").Value = -4862.999899999999
$ws.Range("N138").Value = -60260

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1366.381
$ws.Range("I32").Value = 936.0548
$ws.Range("K32").Value = 936.0548
$ws.Range("M32").Value = -649.0548
# Row 48
$ws.Range("H48").Value = 125684
$ws.Range("J48").Value = 125684
$ws.Range("L48").Value = 125684
$ws.Range("N48").Value = -126452

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 11
$ws.Range("H11").Value = 273
$ws.Range("I11").Value = 99.8
$ws.Range("J11").Value = 561.6667
$ws.Range("K11").Value = 99.8
$ws.Range("L11").Value = 561.6667
$ws.Range("M11").Value = 40.2
$ws.Range("N11").Value = -841.6667
# Row 37
$ws.Range("H37").Value = 865.6
$ws.Range("I37").Value = 566.3333
$ws.Range("J37").Value = 1314.5
$ws.Range("K37").Value = 566.3333
$ws.Range("L37").Value = 1314.5
$ws.Range("M37").Value = -429.3333
$ws.Range("N37").Value = -1588.5
# Row 47
$ws.Range("H47").Value = 175684
$ws.Range("J47").Value = 175684
$ws.Range("L47").Value = 175684
$ws.Range("N47").Value = -176724
# Row 48
$ws.Range("H48").Value = 175684
$ws.Range("J48").Value = 175684
$ws.Range("L48").Value = 175684
$ws.Range("N48").Value = -176514

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 3815
$ws.Range("J6").Value = 4000.5
$ws.Range("L6").Value = 4000.5
$ws.Range("N6").Value = -4226.5
# Row 15
$ws.Range("H15").Value = 5118
$ws.Range("I15").Value = 6157.3335
$ws.Range("J15").Value = 2000
$ws.Range("K15").Value = 6157.3335
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = -5987.3335
$ws.Range("N15").Value = -2340
# Row 59
$ws.Range("H59").Value = 60500
$ws.Range("I59").Value = 21000
$ws.Range("K59").Value = 21000
$ws.Range("M59").Value = -19855
# Row 68
$ws.Range("H68").Value = 175000
$ws.Range("I68").Value = 150000
$ws.Range("J68").Value = 200000
$ws.Range("K68").Value = 150000
$ws.Range("L68").Value = 200000
$ws.Range("M68").Value = -149251
$ws.Range("N68").Value = -201498
# Row 71
$ws.Range("H71").Value = 175000
$ws.Range("I71").Value = 150000
$ws.Range("J71").Value = 200000
$ws.Range("K71").Value = 450000
$ws.Range("L71").Value = 600000
$ws.Range("M71").Value = -446256
$ws.Range("N71").Value = -607488
# Row 87
$ws.Range("H87").Value = 66924.75
$ws.Range("J87").Value = 55233
$ws.Range("L87").Value = 55233
$ws.Range("N87").Value = -57605
# Row 90
$ws.Range("H90").Value = 66924.75
$ws.Range("J90").Value = 55233
$ws.Range("L90").Value = 165699
$ws.Range("N90").Value = -177555
# Row 108
$ws.Range("H108").Value = 60684
$ws.Range("J108").Value = 60684
$ws.Range("L108").Value = 60684
$ws.Range("N108").Value = -68364

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 576.05
$ws.Range("J98").Value = 629.2
$ws.Range("L98").Value = 1887.6
$ws.Range("N98").Value = -4883.6
# Row 107
$ws.Range("H107").Value = 980.8
$ws.Range("I107").Value = 353
$ws.Range("J107").Value = 1399.3334
$ws.Range("K107").Value = 1059
$ws.Range("L107").Value = 4198.0002
$ws.Range("M107").Value = 861
$ws.Range("N107").Value = -8038.0002
# Row 112
$ws.Range("H112").Value = 7616.6665
$ws.Range("I112").Value = 1425
$ws.Range("J112").Value = 20000
$ws.Range("K112").Value = 4275
$ws.Range("L112").Value = 60000
$ws.Range("M112").Value = -3167
$ws.Range("N112").Value = -62216
# Row 116
$ws.Range("H116").Value = 4714.5
$ws.Range("I116").Value = 4714.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 14143.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -10701.5
$ws.Range("N116").ClearContents()
# Row 117
$ws.Range("H117").Value = 803.3333
$ws.Range("I117").Value = 409.6
$ws.Range("J117").Value = 2772
$ws.Range("K117").Value = 1228.8
$ws.Range("L117").Value = 8316
$ws.Range("M117").Value = 2213.2
$ws.Range("N117").Value = -15200
# Row 118
$ws.Range("H118").Value = 100
$ws.Range("I118").Value = 100
$ws.Range("K118").Value = 300
$ws.Range("M118").Value = 943

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 136
$ws.Range("H136").Value = 70990.25
$ws.Range("J136").Value = 70990.25
$ws.Range("L136").Value = 212970.75
$ws.Range("N136").Value = -218070.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7460.3335
$ws.Range("I7").Value = 6818.909
$ws.Range("J7").Value = 9224.25
$ws.Range("K7").Value = 6818.909
$ws.Range("L7").Value = 9224.25
$ws.Range("M7").Value = -6706.909
$ws.Range("N7").Value = -9448.25
# Row 40
$ws.Range("H40").Value = 6524.316
$ws.Range("I40").Value = 6389.6113
$ws.Range("K40").Value = 6389.6113
$ws.Range("M40").Value = -6253.6113
# Row 97
$ws.Range("H97").Value = 23549.5
$ws.Range("J97").Value = 23549.5
$ws.Range("L97").Value = 23549.5
$ws.Range("N97").Value = -25531.5
# Row 104
$ws.Range("H104").Value = 22056.5
$ws.Range("J104").Value = 23777
$ws.Range("L104").Value = 23777
$ws.Range("N104").Value = -30765
# Row 126
$ws.Range("H126").Value = 7460.3335
$ws.Range("I126").Value = 6818.909
$ws.Range("J126").Value = 9224.25
$ws.Range("K126").Value = 20456.727
$ws.Range("L126").Value = 27672.75
$ws.Range("M126").Value = -17986.727
$ws.Range("N126").Value = -32612.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 72885.63
$ws.Range("I2").Value = 85679.38
$ws.Range("J2").Value = 45165.832
$ws.Range("K2").Value = 85679.38
$ws.Range("L2").Value = 45165.832
$ws.Range("M2").Value = -85567.38
$ws.Range("N2").Value = -45389.832
# Row 107
$ws.Range("H107").Value = 1807.96
$ws.Range("I107").Value = 1433.3334
$ws.Range("J107").Value = 2369.9
$ws.Range("K107").Value = 4300.0002
$ws.Range("L107").Value = 7109.700000000001
$ws.Range("M107").Value = -2380.0002
$ws.Range("N107").Value = -10949.7
# Row 122
$ws.Range("H122").Value = 2052.9167
$ws.Range("I122").Value = 1648.8966
$ws.Range("J122").Value = 3726.7144
$ws.Range("K122").Value = 4946.6898
$ws.Range("L122").Value = 11180.1432
$ws.Range("M122").Value = -2496.6898
$ws.Range("N122").Value = -16080.1432
# Row 126
$ws.Range("H126").Value = 4218.0713
$ws.Range("I126").Value = 3175.5
$ws.Range("K126").Value = 9526.5
$ws.Range("M126").Value = -7056.5

Write-Host "Famfrit_Profits update applied"